$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This adds 6 new "Residual Category" rows (Don't Know, Refused to Answer,
# Repeated Value, Response Unidentifiable, Response Outside Scope, Not Stated)
# to the bottom of the Ethnicity NZSCC classification table (rows 182-187),
# following the same Code/Description/Level_*_Code/Level_*_Description layout
# as the rest of the sheet. All "Code" style columns (A,C,E,G) must stay text
# (e.g. "94444", not the number 94444), matching every other row in the sheet.
#
# Simply assigning numeric-looking strings via Range.Value lets Excel convert
# them to real numbers, so instead we type each value into a scratch area as
# a formula ( ="94444" ) which always evaluates to text, then paste-special
# just the *values* into the real destination cells. That keeps the cells as
# proper shared-string/text cells without picking up any new number format.

# Build helper grid with formulas to force text (avoid numeric auto-conversion)
$ws.Range('Z1').Formula = '="94444"'
$ws.Range('AA1').Formula = '="Don''t Know"'
$ws.Range('AB1').Formula = '="944"'
$ws.Range('AC1').Formula = '="Don''t Know"'
$ws.Range('AD1').Formula = '="94"'
$ws.Range('AE1').Formula = '="Don''t Know"'
$ws.Range('AF1').Formula = '="9"'
$ws.Range('AG1').Formula = '="Residual Categories"'
$ws.Range('Z2').Formula = '="95555"'
$ws.Range('AA2').Formula = '="Refused to Answer"'
$ws.Range('AB2').Formula = '="955"'
$ws.Range('AC2').Formula = '="Refused to Answer"'
$ws.Range('AD2').Formula = '="95"'
$ws.Range('AE2').Formula = '="Refused to Answer"'
$ws.Range('AF2').Formula = '="9"'
$ws.Range('AG2').Formula = '="Residual Categories"'
$ws.Range('Z3').Formula = '="96666"'
$ws.Range('AA3').Formula = '="Repeated Value"'
$ws.Range('AB3').Formula = '="966"'
$ws.Range('AC3').Formula = '="Repeated Value"'
$ws.Range('AD3').Formula = '="96"'
$ws.Range('AE3').Formula = '="Repeated Value"'
$ws.Range('AF3').Formula = '="9"'
$ws.Range('AG3').Formula = '="Residual Categories"'
$ws.Range('Z4').Formula = '="97777"'
$ws.Range('AA4').Formula = '="Response Unidentifiable"'
$ws.Range('AB4').Formula = '="977"'
$ws.Range('AC4').Formula = '="Response Unidentifiable"'
$ws.Range('AD4').Formula = '="97"'
$ws.Range('AE4').Formula = '="Response Unidentifiable"'
$ws.Range('AF4').Formula = '="9"'
$ws.Range('AG4').Formula = '="Residual Categories"'
$ws.Range('Z5').Formula = '="98888"'
$ws.Range('AA5').Formula = '="Response Outside Scope"'
$ws.Range('AB5').Formula = '="988"'
$ws.Range('AC5').Formula = '="Response Outside Scope"'
$ws.Range('AD5').Formula = '="98"'
$ws.Range('AE5').Formula = '="Response Outside Scope"'
$ws.Range('AF5').Formula = '="9"'
$ws.Range('AG5').Formula = '="Residual Categories"'
$ws.Range('Z6').Formula = '="99999"'
$ws.Range('AA6').Formula = '="Not Stated"'
$ws.Range('AB6').Formula = '="999"'
$ws.Range('AC6').Formula = '="Not Stated"'
$ws.Range('AD6').Formula = '="99"'
$ws.Range('AE6').Formula = '="Not Stated"'
$ws.Range('AF6').Formula = '="9"'
$ws.Range('AG6').Formula = '="Residual Categories"'

# Copy helper grid and paste as values into target rows 182-187, columns A-H
$ws.Range('Z1:AG6').Copy() | Out-Null
$ws.Range('A182:H187').PasteSpecial(-4163) | Out-Null  # xlPasteValues

# Clear helper grid
$ws.Range('Z1:AG6').ClearContents() | Out-Null

# Column I (Definition) is blank text (an empty shared string, same as every
# other row). A plain empty Value is treated as a truly blank cell, so force
# a text cell via the leading-apostrophe trick, then re-apply I181's (plain,
# unstyled) format on top so the new cells don't end up with a different
# style than the rest of the sheet.
$ws.Range('I182:I187').Value = "'"
$ws.Range('I181').Copy() | Out-Null
$ws.Range('I182:I187').PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false

